$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.347.19"
$ws.Range("E2").Value = "'  -0.14%  "
$ws.Range("D3").Value = "'1.841.80"
$ws.Range("E3").Value = "'  -0.23%  "
$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'239.68"
$ws.Range("E5").Value = "'  -0.45%  "
$ws.Range("D6").Value = "'0.6291"
$ws.Range("E6").Value = "'  -0.14%  "
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("D8").Value = "'0.07437"
$ws.Range("E8").Value = "'  -0.69%  "
$ws.Range("E9").Value = "'  -0.54%  "
$ws.Range("D10").Value = "'24.91"
$ws.Range("E10").Value = "'  +2.10%  "
$ws.Range("D11").Value = "'0.07730"
$ws.Range("E11").Value = "'  +0.01%  "
$ws.Range("D12").Value = "'1.844.04"
$ws.Range("E12").Value = "'  -0.13%  "
$ws.Range("D13").Value = "'4.964"
$ws.Range("E13").Value = "'  -0.80%  "
$ws.Range("D14").Value = "'0.6756"
$ws.Range("E14").Value = "'  -0.46%  "
$ws.Range("D15").Value = "'0.00001024"
$ws.Range("E15").Value = "'  +0.32%  "
$ws.Range("D16").Value = "'81.50"
$ws.Range("E16").Value = "'  -0.76%  "
$ws.Range("D17").Value = "'6.243"
$ws.Range("E17").Value = "'  +1.62%  "
$ws.Range("D18").Value = "'29.382.92"
$ws.Range("E18").Value = "'  -0.14%  "
$ws.Range("D19").Value = "'228.86"
$ws.Range("E19").Value = "'  +0.22%  "
$ws.Range("E20").Value = "'  -0.19%  "
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "'  +0.01%  "
$ws.Range("D22").Value = "'7.355"
$ws.Range("E23").Value = "'  +0.08%  "
$ws.Range("E24").Value = "'  -0.64%  "
$ws.Range("D25").Value = "'8.480"
$ws.Range("E25").Value = "'  +0.70%  "
$ws.Range("E26").Value = "'  -1.97%  "
$ws.Range("D27").Value = "'17.41"
$ws.Range("E27").Value = "'  -0.78%  "
$ws.Range("D28").Value = "'0.07068"
$ws.Range("E28").Value = "'  +12.48%  "
$ws.Range("E29").Value = "'  +5.70%  "
$ws.Range("D30").Value = "'1.480"
$ws.Range("E30").Value = "'  +0.41%  "
$ws.Range("D31").Value = "'4.051"
$ws.Range("E31").Value = "'  -1.03%  "
$ws.Range("D32").Value = "'4.030"
$ws.Range("E32").Value = "'  -0.70%  "
$ws.Range("D33").Value = "'1.825"
$ws.Range("E33").Value = "'  +0.29%  "
$ws.Range("E34").Value = "'  -0.26%  "
$ws.Range("D35").Value = "'0.6974"
$ws.Range("E35").Value = "'  +0.08%  "
$ws.Range("D36").Value = "'2.582"
$ws.Range("E36").Value = "'  +0.00%  "
$ws.Range("D37").Value = "'0.01841"
$ws.Range("E37").Value = "'  +1.06%  "
$ws.Range("E38").Value = "'  -0.64%  "
$ws.Range("D39").Value = "'1.235.17"
$ws.Range("E39").Value = "'  -1.71%  "
$ws.Range("D40").Value = "'6.803"
$ws.Range("E40").Value = "'  +3.82%  "
$ws.Range("D41").Value = "'0.9278"
$ws.Range("E41").Value = "'  +2.04%  "
$ws.Range("D42").Value = "'0.9999"
$ws.Range("E42").Value = "'  +0.01%  "
$ws.Range("D43").Value = "'1.996.62"
$ws.Range("E43").Value = "'  -0.60%  "
$ws.Range("D44").Value = "'100.77"
$ws.Range("E44").Value = "'  -0.63%  "
$ws.Range("D45").Value = "'65.28"
$ws.Range("E45").Value = "'  -1.60%  "
$ws.Range("E46").Value = "'  +2.01%  "
$ws.Range("D47").Value = "'7.017"
$ws.Range("E47").Value = "'  -0.45%  "
$ws.Range("E48").Value = "'  +1.22%  "
$ws.Range("E49").Value = "'  -2.90%  "
$ws.Range("D50").Value = "'8.898"
$ws.Range("E50").Value = "'  -1.63%  "
$ws.Range("D51").Value = "'0.3914"
$ws.Range("E51").Value = "'  -0.67%  "
